# Added main, modified test, and added function for postfix eval and cal
#
# The worksheet gains a new title row above the existing table: a merged
# A1:C1 cell containing "Infix Expression: a*b/(c-a)+d*e", where the
# "Infix Expression" part is bold and the rest is regular. Every
# pre-existing row shifts down by one. The former header row (now row 2)
# becomes bold as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push all existing rows down by inserting a new row 1.
$ws.Rows("1:1").Insert()

# New merged title cell with the infix-expression caption.
$ws.Range("A1:C1").Value = "Infix Expression: a*b/(c-a)+d*e"
$ws.Range("A1:C1").HorizontalAlignment = -4108
$ws.Range("A1:C1").Merge()

# Bold just the "Infix Expression" label (first 16 characters).
$ws.Range("A1").Characters(1, 16).Font.Bold = $true

# The old header row (Next Character.../Postfix Form/Operator Stack...)
# is now row 2; make it bold.
$ws.Range("A2:C2").Font.Bold = $true

# Best-effort: touch page setup (original workbook picked up printer
# defaults here; orientation is the only part we can influence headlessly).
$ws.PageSetup.Orientation = 1

# Restore the cursor/selection to where the author left it.
$ws.Range("C8").Select() | Out-Null
